$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "heaven" sound effect (row 12) used a CC BY-NC (non-commercial) licensed
# clip from Robinhood76. Replace it with a CC0 clip from random_intruder.
# Set E12 (source URL) before D12 (copyrighter name) so the new shared
# strings are appended to the sharedStrings table in the same order as the
# target workbook (URL first, then name).
$ws.Range("E12").Value = "https://freesound.org/people/random_intruder/sounds/392172/"
$ws.Range("D12").Value = "random_intruder"
$ws.Range("F12").Value = "CC0 1.0(No Copyright, Public Domain Dedication)/ https://creativecommons.org/publicdomain/zero/1.0/"

# Update the saved cursor/selection position to match the author's final
# selection when they saved the workbook.
$ws.Range("F21").Select()
